$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("recharge").Range("I2").Value = '{''code'': ''10001'', ''data'': {''mobilephone'': ''17751810000'', ''type'': ''1'', ''regtime'': ''2019-01-21 10:35:36.0'', ''pwd'': ''E10ADC3949BA59ABBE56E057F20F883E'', ''id'': 1115516, ''regname'': ''小蜜蜂'', ''leaveamount'': ''79553.00''}, ''msg'': ''充值成功'', ''status'': 1}'
$wb.Worksheets.Item("recharge").Range("I3").Value = '{''code'': ''10001'', ''data'': {''mobilephone'': ''17751810000'', ''type'': ''1'', ''regtime'': ''2019-01-21 10:35:36.0'', ''pwd'': ''E10ADC3949BA59ABBE56E057F20F883E'', ''id'': 1115516, ''regname'': ''小蜜蜂'', ''leaveamount'': ''80054.00''}, ''msg'': ''充值成功'', ''status'': 1}'
$wb.Worksheets.Item("recharge").Range("I6").Value = '{''code'': ''10001'', ''data'': {''mobilephone'': ''17751810001'', ''type'': ''1'', ''regtime'': ''2019-01-21 10:30:31.0'', ''pwd'': ''F1887D3F9E6EE7A32FE5E76F4AB80D63'', ''id'': 1115509, ''regname'': ''å¤\x9cé\x9b¨å£°ç\x83¦'', ''leaveamount'': ''52920.00''}, ''msg'': ''充值成功'', ''status'': 1}'
$wb.Worksheets.Item("recharge").Range("I12").Value = '{''code'': ''10001'', ''data'': {''mobilephone'': ''17751810000'', ''type'': ''1'', ''regtime'': ''2019-01-21 10:35:36.0'', ''pwd'': ''E10ADC3949BA59ABBE56E057F20F883E'', ''id'': 1115516, ''regname'': ''小蜜蜂'', ''leaveamount'': ''80554.00''}, ''msg'': ''充值成功'', ''status'': 1}'
$wb.Worksheets.Item("recharge").Range("I13").Value = '{''code'': ''10001'', ''data'': {''mobilephone'': ''17751810000'', ''type'': ''1'', ''regtime'': ''2019-01-21 10:35:36.0'', ''pwd'': ''E10ADC3949BA59ABBE56E057F20F883E'', ''id'': 1115516, ''regname'': ''小蜜蜂'', ''leaveamount'': ''81055.00''}, ''msg'': ''充值成功'', ''status'': 1}'
$wb.Worksheets.Item("recharge").Range("I16").Value = '{''code'': ''10001'', ''data'': {''mobilephone'': ''17751810001'', ''type'': ''1'', ''regtime'': ''2019-01-21 10:30:31.0'', ''pwd'': ''F1887D3F9E6EE7A32FE5E76F4AB80D63'', ''id'': 1115509, ''regname'': ''å¤\x9cé\x9b¨å£°ç\x83¦'', ''leaveamount'': ''53424.00''}, ''msg'': ''充值成功'', ''status'': 1}'
$wb.Worksheets.Item("withdraw").Range("I3").Value = '{''code'': ''10001'', ''data'': {''mobilephone'': ''17751810779'', ''type'': ''1'', ''regtime'': ''2019-01-21 18:01:06.0'', ''pwd'': ''E10ADC3949BA59ABBE56E057F20F883E'', ''id'': 1115661, ''regname'': ''小蜜蜂'', ''leaveamount'': ''499735.00''}, ''msg'': ''取现成功'', ''status'': 1}'
$wb.Worksheets.Item("withdraw").Range("I14").Value = '{''code'': ''10001'', ''data'': {''mobilephone'': ''17751810779'', ''type'': ''1'', ''regtime'': ''2019-01-21 18:01:06.0'', ''pwd'': ''E10ADC3949BA59ABBE56E057F20F883E'', ''id'': 1115661, ''regname'': ''小蜜蜂'', ''leaveamount'': ''499730.00''}, ''msg'': ''取现成功'', ''status'': 1}'
